
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the "Week 3" column (N) ---
# Header cell N2: copy M2's format (bold/centered/filled header style) then set the text.
$ws.Range("M2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("N2").Value = 'Week 3'

# Body cells N3:N5: copy the wrap-text body style used by the Week-2 column, then
# fill in the Week-3 progress notes for teams 1-3.
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N3").Value = 'The team is using 10% over their dataset for now. They are using oversampling technique to handle the class imbalance issue. They have used one-hot encoder for categorical data. They have tried Decision Tree and Naive Bayes until now. They are working on more models. Since there is high correlation between some columns, the team is planning to test dimensionality reduction.'

$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N4").Value = 'The team has performed some more EDA to try to separate the distribution of votes per party.'

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N5").Value = 'After separating the “dow”, “cities”, and “genres” columns, the team got a very large number of features. The team realized that the total time is the same for “dow”, “cities”, and “genres”. So, the team decided to go only with the “dow” column. They have used Logistic Regression but it doesn''t perform well. So, they are trying out different feature engineering techniques now.'

$excel.CutCopyMode = 0

# --- Team 7 (row 9) didn't have a Week-2 note yet; add it now ---
$ws.Range("M9").Value = 'Since the previous dataset didn’t have enough information to be used in a Data Science project, the team has obtained a new dataset about NYC AirBnB prices. They are currently preprocessing the dataset and will move on to EDA after this part.'
$ws.Rows.Item(9).RowHeight = 85.05

# --- Column widths (Week 2 / Week 3 columns) ---
$ws.Columns.Item(13).ColumnWidth = 31.43
$ws.Columns.Item(14).ColumnWidth = 31.43

# --- View state: scroll so column K is at the left edge and select N4 ---
$win = $excel.Windows.Item(1)
$win.ScrollRow = 3
$win.ScrollColumn = 11
$ws.Range("N4").Select()
